$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell values that changed between the old and new clustering run
$ws.Range("B2").Value = 0.02186588921282799
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0.002695417789757413
$ws.Range("G2").Value = 0.002989130434782609
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.01412776412776413
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("B3").Value = 0.03790087463556849
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 0.008984725965858042
$ws.Range("G3").Value = 0.02010869565217389
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.01228501228501228
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("B4").Value = 0.08527696793002894
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.02361142343152681
$ws.Range("F4").Value = 0.002695417789757413
$ws.Range("G4").Value = 0.02880434782608689
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.006756756756756756
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.00830367734282325
$ws.Range("B5").Value = 0.02551020408163266
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.1425680233865518
$ws.Range("F5").Value = 0.0008984725965858042
$ws.Range("G5").Value = 0.03695652173913032
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.002457002457002457
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.2253855278766305
$ws.Range("B6").Value = 0.009475218658892131
$ws.Range("D6").Value = 0
$ws.Range("F6").Value = 0.001796945193171608
$ws.Range("G6").Value = 0.005978260869565219
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0.004914004914004914
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("C7").Value = 0.00759493670886076
$ws.Range("F7").Value = 0.008984725965858042
$ws.Range("H7").Value = 0.005747126436781609
$ws.Range("I7").Value = 0.01781326781326782
$ws.Range("B8").Value = 0.167638483965014
$ws.Range("C8").Value = 0.02025316455696203
$ws.Range("D8").Value = 0.1070384528895878
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0.02246181491464511
$ws.Range("G8").Value = 0.372554347826085
$ws.Range("H8").Value = 0.02155172413793103
$ws.Range("I8").Value = 0.0380835380835381
$ws.Range("J8").Value = 0.01814516129032258
$ws.Range("K8").Value = 0.01304863582443653
$ws.Range("F9").Value = 0.004492362982929021
$ws.Range("H9").Value = 0.001436781609195402
$ws.Range("I9").Value = 0.0208845208845209
$ws.Range("B10").Value = 0.07361516034985406
$ws.Range("C10").Value = 0.03291139240506329
$ws.Range("D10").Value = 0.006521250281088375
$ws.Range("E10").Value = 0.005571030640668524
$ws.Range("F10").Value = 0.04222821203953284
$ws.Range("G10").Value = 0.01874999999999999
$ws.Range("H10").Value = 0.03735632183908043
$ws.Range("I10").Value = 0.1437346437346439
$ws.Range("J10").Value = 0.03629032258064514
$ws.Range("K10").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0.03544303797468355
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0.02425876010781673
$ws.Range("H12").Value = 0.03735632183908043
$ws.Range("J12").Value = 0.008064516129032258
$ws.Range("B13").Value = 0.002915451895043732
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0.01416685405891613
$ws.Range("G13").Value = 0.005434782608695654
$ws.Range("I13").Value = 0.0006142506142506142
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("B14").Value = 0.01676384839650146
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0.03485495839892048
$ws.Range("F14").Value = 0.001796945193171608
$ws.Range("G14").Value = 0.00706521739130435
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0.004914004914004914
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0.2550415183867135
$ws.Range("D15").Value = 0
$ws.Range("I15").Value = 0.03501228501228504
$ws.Range("B16").Value = 0.05903790087463546
$ws.Range("C16").Value = 0.01265822784810127
$ws.Range("D16").Value = 0.03665392399370346
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0.01796945193171608
$ws.Range("G16").Value = 0.02853260869565211
$ws.Range("H16").Value = 0.01005747126436782
$ws.Range("I16").Value = 0.02825552825552829
$ws.Range("J16").Value = 0.006048387096774193
$ws.Range("K16").Value = 0
$ws.Range("B17").Value = 0.003644314868804665
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = 0.0008984725965858042
$ws.Range("G17").Value = 0.0008152173913043479
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0.005063291139240506
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0.001796945193171608
$ws.Range("H18").Value = 0.002873563218390805
$ws.Range("J18").Value = 0.002016129032258064
$ws.Range("B19").Value = 0.01749271137026239
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0.01709017315043849
$ws.Range("F19").Value = 0.0008984725965858042
$ws.Range("G19").Value = 0.01630434782608696
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0.004914004914004914
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0.01245551601423487
$ws.Range("B20").Value = 0.008017492711370264
$ws.Range("D20").Value = 0
$ws.Range("G20").Value = 0.005978260869565219
$ws.Range("I20").Value = 0.004914004914004914
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("B21").Value = 0.007288629737609331
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0.009894310771306499
$ws.Range("F21").Value = 0.001796945193171608
$ws.Range("G21").Value = 0.005706521739130437
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0.003071253071253071
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("C22").Value = 0.02784810126582278
$ws.Range("E22").Value = 0.005571030640668524
$ws.Range("F22").Value = 0.04582210242587607
$ws.Range("H22").Value = 0.06034482758620684
$ws.Range("I22").Value = 0.001228501228501228
$ws.Range("J22").Value = 0.01209677419354839
$ws.Range("B23").Value = 0.05903790087463546
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("F23").Value = 0.02695417789757414
$ws.Range("G23").Value = 0.03614130434782598
$ws.Range("H23").Value = 0.002873563218390805
$ws.Range("I23").Value = 0.02702702702702706
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0.005390835579514825
$ws.Range("H24").Value = 0.005747126436781609
$ws.Range("I24").Value = 0.0006142506142506142
$ws.Range("J24").Value = 0.004032258064516129
$ws.Range("B32").Value = 0.06268221574344011
$ws.Range("C32").Value = 0.03544303797468355
$ws.Range("D32").Value = 0.02675961322239704
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0.07816711590296507
$ws.Range("G32").Value = 0.01086956521739131
$ws.Range("H32").Value = 0.05028735632183904
$ws.Range("I32").Value = 0.08968058968058985
$ws.Range("J32").Value = 0.008064516129032258
$ws.Range("K32").Value = 0.08540925266903901
$ws.Range("B33").Value = 0.02040816326530613
$ws.Range("C33").Value = 0.1341772151898735
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0.1644204851752018
$ws.Range("H33").Value = 0.1623563218390806
$ws.Range("I33").Value = 0.1246928746928752
$ws.Range("J33").Value = 0.07056451612903221
$ws.Range("K33").Value = 0
$ws.Range("B34").Value = 0.002186588921282799
$ws.Range("C34").Value = 0
$ws.Range("G34").Value = 0.0008152173913043479
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0

# The new clustering run dropped the "joint regime area" rows (36-40);
# remove them so the sheet's used range shrinks to A1:K35
$ws.Rows("36:40").Delete()
